$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from SCD0330 to SCD0024
$ws.Name = "SCD0024"

# Update the TC_ID cell (B2) from "DGS-345" to "SCD0024-009"
$ws.Range("B2").Value = "SCD0024-009"

# Widen column B to fit the new (wider) content, matching the bestFit width
$ws.Columns("B").ColumnWidth = 14.6

# Move the active selection to B3 (cursor position after edit)
$ws.Range("B3").Select()
